# Update Name of Algo
# Applies the cell value changes described by the source diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B3"  = 6.845000000000001
    "B4"  = 6.612
    "C6"  = -12.52
    "B7"  = 6.411
    "C7"  = -12.918
    "B8"  = 5.823
    "C8"  = -12.241
    "A11" = -21.643
    "A12" = -21.456
    "B12" = 6.694999999999999
    "B14" = 6.567
    "A15" = -21.192
    "C19" = -12.684
    "C21" = -12.734
    "B22" = 6.686
    "C24" = -12.177
    "C25" = -12.69
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
